$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update author name and price for row 2 (identification No 1)
$ws.Range("C2").Value = "Aaron"
$ws.Range("D2").Value = 15

# Update author name for row 3 (identification No 2)
$ws.Range("C3").Value = "Petrosky"
